$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.462.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.175.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.01%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.16%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.175.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.21"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.693.80"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.170.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.442.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.89"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.103"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.05%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.17"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0692"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "412.60"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.871.52"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.60"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.08%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.64"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.58%  "
